# Anaconda distribution.pptx — slide 10 ("Thank you for Listening" closer)
#
# Turns the lone full-bleed picture on slide 10 into a group containing that
# same picture (nudged up/right) plus a new yellow caption rectangle reading
# "Thank you for Listening!", matching the authored edit.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# ---- the pre-existing picture -------------------------------------------------
$pic = $s.Shapes.Item(1)
$pic.Name = "Picture 3"

# ---- new caption rectangle ----------------------------------------------------
$rectLeft   = 4452257 / $EMU_PER_PT
$rectTop    = 3442138 / $EMU_PER_PT
$rectWidth  = 3145972 / $EMU_PER_PT
$rectHeight = 359229  / $EMU_PER_PT

$rect = $s.Shapes.AddShape(1, $rectLeft, $rectTop, $rectWidth, $rectHeight)
$rect.Name = "Rectangle 5"

$rect.Fill.ForeColor.RGB = 6746111   # 0xFFEF66 (R=255,G=239,B=102) -> R+G*256+B*65536
$rect.Line.Visible = $false

$tf = $rect.TextFrame
$tf.TextRange.Text = "Thank you for Listening!"
$tf.VerticalAnchor = 3               # msoAnchorMiddle

$tr = $tf.TextRange
$tr.Font.Name = "Arial"
$tr.Font.NameComplexScript = "Arial"
$tr.Font.Size = 20
$tr.Font.Bold = $true
$tr.Font.Color.RGB = 0               # black == theme tx1 in this deck
$tr.ParagraphFormat.Alignment = 2    # ppAlignCenter

# ---- group the picture with the new rectangle ---------------------------------
$range = $s.Shapes.Range(@($pic.Name, $rect.Name))
$grp = $range.Group()
$grp.Name = "Group 6"

# ---- reposition the picture inside the group, then the group itself -----------
$picInGroup = $grp.GroupItems.Item(1)
$picInGroup.Left = 2468241 / $EMU_PER_PT
$picInGroup.Top  = 350778  / $EMU_PER_PT

$grp.Left = 2654890 / $EMU_PER_PT
$grp.Top  = 1015120 / $EMU_PER_PT
